# Update the "想去人数" (want-to-go count) column F on the relevant sheets
# to match the newly scraped totals.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 205
$ws1.Range("F4").Value = 508
$ws1.Range("F5").Value = 505
$ws1.Range("F7").Value = 2583
$ws1.Range("F8").Value = 441
$ws1.Range("F9").Value = 7028
$ws1.Range("F11").Value = 441
$ws1.Range("F13").Value = 91

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 13
$ws2.Range("F3").Value = 17

# --- Sheet "全部类型" (all types, combined listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 205
$ws4.Range("F4").Value = 508
$ws4.Range("F5").Value = 505
$ws4.Range("F7").Value = 13
$ws4.Range("F8").Value = 17
$ws4.Range("F9").Value = 2583
$ws4.Range("F10").Value = 441
$ws4.Range("F11").Value = 7028
$ws4.Range("F13").Value = 441
$ws4.Range("F17").Value = 91
